# Auto-generated edit script: updates profit-calculation columns (H-N)
# on several rows across all 8 sheets, per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3782.2222
$ws.Range("I62").Value = 3934.2856
$ws.Range("J62").Value = 3250
$ws.Range("K62").Value = 3934.2856
$ws.Range("L62").Value = 3250
$ws.Range("M62").Value = -3310.2856
$ws.Range("N62").Value = -4498

$ws.Range("H65").Value = 3782.2222
$ws.Range("I65").Value = 3934.2856
$ws.Range("J65").Value = 3250
$ws.Range("K65").Value = 19671.428
$ws.Range("L65").Value = 16250
$ws.Range("M65").Value = -16551.428
$ws.Range("N65").Value = -22490

$ws.Range("H116").Value = 2318.2222
$ws.Range("I116").Value = 2010.6666
$ws.Range("J116").Value = 2933.3333
$ws.Range("K116").Value = 2010.6666
$ws.Range("L116").Value = 2933.3333
$ws.Range("M116").Value = 1431.3334
$ws.Range("N116").Value = -9817.3333

$ws.Range("H137").Value = 947.55554
$ws.Range("I137").Value = 905.2414
$ws.Range("K137").Value = 2715.7242
$ws.Range("M137").Value = -165.7242000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 21994.75
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 21994.75
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 21994.75
$ws.Range("M44").ClearContents() | Out-Null
$ws.Range("N44").Value = -22970.75

$ws.Range("H54").Value = 9800
$ws.Range("J54").Value = 9800
$ws.Range("L54").Value = 9800
$ws.Range("N54").Value = -11338

$ws.Range("H61").Value = 2343.7144
$ws.Range("I61").Value = 1729.8572
$ws.Range("K61").Value = 1729.8572
$ws.Range("M61").Value = -1517.8572

$ws.Range("H74").Value = 414.15384
$ws.Range("I74").Value = 358.4
$ws.Range("J74").Value = 600
$ws.Range("K74").Value = 358.4
$ws.Range("L74").Value = 600
$ws.Range("M74").Value = 515.6
$ws.Range("N74").Value = -2348

$ws.Range("H77").Value = 414.15384
$ws.Range("I77").Value = 358.4
$ws.Range("J77").Value = 600
$ws.Range("K77").Value = 1792
$ws.Range("L77").Value = 3000
$ws.Range("M77").Value = 2576
$ws.Range("N77").Value = -11736

$ws.Range("H136").Value = 2343.7144
$ws.Range("I136").Value = 1729.8572
$ws.Range("K136").Value = 5189.571599999999
$ws.Range("M136").Value = -2639.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 48162.547
$ws.Range("I134").Value = 69171.734
$ws.Range("J134").Value = 3142.8572
$ws.Range("K134").Value = 207515.202
$ws.Range("L134").Value = 9428.5716
$ws.Range("M134").Value = -204980.202
$ws.Range("N134").Value = -14498.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1685061.6
$ws.Range("I31").Value = 1682.6945
$ws.Range("K31").Value = 1682.6945
$ws.Range("M31").Value = -1387.6945

$ws.Range("H34").Value = 1685061.6
$ws.Range("I34").Value = 1682.6945
$ws.Range("K34").Value = 1682.6945
$ws.Range("M34").Value = -1480.6945

$ws.Range("H58").Value = 1093.8572
$ws.Range("I58").Value = 1108.5385
$ws.Range("J58").Value = 1070
$ws.Range("K58").Value = 1108.5385
$ws.Range("L58").Value = 1070
$ws.Range("M58").Value = -905.5385000000001
$ws.Range("N58").Value = -1476

$ws.Range("H94").Value = 3425
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 3425
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 3425
$ws.Range("M94").ClearContents() | Out-Null
$ws.Range("N94").Value = -4327

$ws.Range("H105").Value = 1653.5
$ws.Range("I105").Value = 3010
$ws.Range("J105").Value = 1382.2
$ws.Range("K105").Value = 3010
$ws.Range("L105").Value = 1382.2
$ws.Range("M105").Value = -1263
$ws.Range("N105").Value = -4876.2

$ws.Range("H132").Value = 2236
$ws.Range("I132").Value = 1420.9231
$ws.Range("J132").Value = 3199.2727
$ws.Range("K132").Value = 4262.7693
$ws.Range("L132").Value = 9597.8181
$ws.Range("M132").Value = -1732.7693
$ws.Range("N132").Value = -14657.8181

$ws.Range("H134").Value = 966
$ws.Range("I134").Value = 728
$ws.Range("K134").Value = 2184
$ws.Range("M134").Value = 351

$ws.Range("H136").Value = 1093.8572
$ws.Range("I136").Value = 1108.5385
$ws.Range("J136").Value = 1070
$ws.Range("K136").Value = 3325.6155
$ws.Range("L136").Value = 3210
$ws.Range("M136").Value = -775.6155000000003
$ws.Range("N136").Value = -8310

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1001.9388
$ws.Range("I68").Value = 840
$ws.Range("J68").Value = 1406.7858
$ws.Range("K68").Value = 2520
$ws.Range("L68").Value = 4220.357400000001
$ws.Range("M68").Value = -1709
$ws.Range("N68").Value = -5842.357400000001

$ws.Range("H71").Value = 1001.9388
$ws.Range("I71").Value = 840
$ws.Range("J71").Value = 1406.7858
$ws.Range("K71").Value = 7560
$ws.Range("L71").Value = 12661.0722
$ws.Range("M71").Value = -3504
$ws.Range("N71").Value = -20773.0722

$ws.Range("H107").Value = 456.10715
$ws.Range("J107").Value = 560.4667
$ws.Range("L107").Value = 1681.4001
$ws.Range("N107").Value = -5521.4001

$ws.Range("H137").Value = 27793252
$ws.Range("I137").Value = 691.53845
$ws.Range("J137").Value = 100053910
$ws.Range("K137").Value = 2074.61535
$ws.Range("L137").Value = 300161730
$ws.Range("M137").Value = 3025.38465
$ws.Range("N137").Value = -300171930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 35000
$ws.Range("J21").Value = 35000
$ws.Range("L21").Value = 35000
$ws.Range("N21").Value = -35346

$ws.Range("H30").Value = 35000
$ws.Range("J30").Value = 35000
$ws.Range("L30").Value = 35000
$ws.Range("N30").Value = -35210

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents() | Out-Null

$ws.Range("H55").Value = 930
$ws.Range("I55").Value = 930
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 930
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -603
$ws.Range("N55").ClearContents() | Out-Null

$ws.Range("H126").Value = 10057.6
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 11322
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 33966
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -38906

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 16000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 16000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 16000
$ws.Range("M45").ClearContents() | Out-Null
$ws.Range("N45").Value = -16814

$ws.Range("H46").Value = 3612.375
$ws.Range("I46").Value = 5499.5
$ws.Range("J46").Value = 2983.3333
$ws.Range("K46").Value = 5499.5
$ws.Range("L46").Value = 2983.3333
$ws.Range("M46").Value = -5311.5
$ws.Range("N46").Value = -3359.3333

$ws.Range("H132").Value = 9506.73
$ws.Range("I132").Value = 13125.6875
$ws.Range("J132").Value = 3716.4
$ws.Range("K132").Value = 39377.0625
$ws.Range("L132").Value = 11149.2
$ws.Range("M132").Value = -36847.0625
$ws.Range("N132").Value = -16209.2

$ws.Range("H136").Value = 6635.5
$ws.Range("I136").Value = 13560.25
$ws.Range("J136").Value = 2678.5
$ws.Range("K136").Value = 40680.75
$ws.Range("L136").Value = 8035.5
$ws.Range("M136").Value = -38130.75
$ws.Range("N136").Value = -13135.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1458.6
$ws.Range("I132").Value = 962.46155
$ws.Range("J132").Value = 1996.0834
$ws.Range("K132").Value = 2887.38465
$ws.Range("L132").Value = 5988.2502
$ws.Range("M132").Value = -357.38465
$ws.Range("N132").Value = -11048.2502

$ws.Range("H136").Value = 3038.5454
$ws.Range("I136").Value = 5386.8
$ws.Range("J136").Value = 1081.6666
$ws.Range("K136").Value = 16160.4
$ws.Range("L136").Value = 3244.9998
$ws.Range("M136").Value = -13610.4
$ws.Range("N136").Value = -8344.9998
